# Add "NA" values under the duplicate_image_filename column (column E)
# for every data row (rows 2 through 21) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2..21) {
    $ws.Cells.Item($r, 5).Value = "NA"
}
